$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (Carou): position changes from "Meio-Campo" to "Zagueira"
$ws.Range("B15").Value = "Zagueira"

# New row 20: Helen, Atacante, 1
$ws.Range("A20").Value = "Helen"
$ws.Range("B20").Value = "Atacante"
$ws.Range("C20").Value = 1

# New row 21: Isadora, Atacante, 1
$ws.Range("A21").Value = "Isadora"
$ws.Range("B21").Value = "Atacante"
$ws.Range("C21").Value = 1

# Apply the same style (border) as the rest of the data rows to the new rows
$ws.Range("A19:C19").Copy()
$ws.Range("A20:C21").PasteSpecial(-4122)

# Update the active cell selection
$ws.Range("A8").Select()
